$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1544401544401544
$ws.Range("C2").Value = 0.6177606177606177
$ws.Range("J2").Value = 0.007722007722007722
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.07722007722007722
$ws.Range("C3").Value = 0.06358381502890173
$ws.Range("J3").Value = 0.04046242774566474
$ws.Range("P3").Value = 0.7687861271676301
$ws.Range("S3").Value = 0.1271676300578035
$ws.Range("J4").Value = 0.08
$ws.Range("P4").Value = 0.68
$ws.Range("S4").Value = 0.24
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.04591836734693878
$ws.Range("D6").Value = 0.01530612244897959
$ws.Range("F6").Value = 0.05102040816326531
$ws.Range("J6").Value = 0.2244897959183673
$ws.Range("O6").Value = 0.01020408163265306
$ws.Range("Q6").Value = 0.1479591836734694
$ws.Range("R6").Value = 0.09693877551020408
$ws.Range("S6").Value = 0.4081632653061225
$ws.Range("B7").Value = 0.1525423728813559
$ws.Range("D7").Value = 0.01129943502824859
$ws.Range("F7").Value = 0.03954802259887006
$ws.Range("J7").Value = 0.1073446327683616
$ws.Range("O7").Value = 0.02824858757062147
$ws.Range("Q7").Value = 0.1468926553672316
$ws.Range("R7").Value = 0.1186440677966102
$ws.Range("S7").Value = 0.3954802259887006
$ws.Range("B8").Value = 0.1098901098901099
$ws.Range("D8").Value = 0.01098901098901099
$ws.Range("F8").Value = 0.08516483516483517
$ws.Range("J8").Value = 0.1043956043956044
$ws.Range("O8").Value = 0.01923076923076923
$ws.Range("Q8").Value = 0.1538461538461539
$ws.Range("R8").Value = 0.1016483516483516
$ws.Range("S8").Value = 0.4148351648351649
$ws.Range("B9").Value = 0.072992700729927
$ws.Range("D9").Value = 0.0072992700729927
$ws.Range("F9").Value = 0.0948905109489051
$ws.Range("J9").Value = 0.1021897810218978
$ws.Range("O9").Value = 0.0072992700729927
$ws.Range("Q9").Value = 0.1605839416058394
$ws.Range("R9").Value = 0.1605839416058394
$ws.Range("S9").Value = 0.3941605839416059
$ws.Range("B10").Value = 0.1202938475665748
$ws.Range("D10").Value = 0.01469237832874196
$ws.Range("E10").Value = 0.002754820936639119
$ws.Range("F10").Value = 0.08999081726354453
$ws.Range("J10").Value = 0.1019283746556474
$ws.Range("O10").Value = 0.009182736455463728
$ws.Range("Q10").Value = 0.1891643709825528
$ws.Range("R10").Value = 0.09182736455463728
$ws.Range("S10").Value = 0.3801652892561984
$ws.Range("G11").Value = 0.1520912547528517
$ws.Range("J11").Value = 0.07224334600760456
$ws.Range("K11").Value = 0.1863117870722434
$ws.Range("L11").Value = 0.5817490494296578
$ws.Range("S11").Value = 0.007604562737642586
$ws.Range("G12").Value = 0.7806451612903226
$ws.Range("J12").Value = 0.167741935483871
$ws.Range("K12").Value = 0.01935483870967742
$ws.Range("L12").Value = 0.01935483870967742
$ws.Range("S12").Value = 0.01290322580645161
$ws.Range("G13").Value = 0.5294117647058824
$ws.Range("J13").Value = 0.3823529411764706
$ws.Range("S13").Value = 0.08823529411764706
$ws.Range("F15").Value = 0.006134969325153374
$ws.Range("H15").Value = 0.2085889570552147
$ws.Range("I15").Value = 0.03680981595092025
$ws.Range("J15").Value = 0.3619631901840491
$ws.Range("K15").Value = 0.09202453987730061
$ws.Range("O15").Value = 0.05521472392638037
$ws.Range("S15").Value = 0.2392638036809816
$ws.Range("F16").Value = 0.01648351648351648
$ws.Range("H16").Value = 0.1318681318681319
$ws.Range("I16").Value = 0.06043956043956044
$ws.Range("J16").Value = 0.4010989010989011
$ws.Range("K16").Value = 0.1428571428571428
$ws.Range("M16").Value = 0.03846153846153846
$ws.Range("O16").Value = 0.05494505494505494
$ws.Range("S16").Value = 0.1538461538461539
$ws.Range("F17").Value = 0.008849557522123894
$ws.Range("H17").Value = 0.1769911504424779
$ws.Range("I17").Value = 0.08554572271386431
$ws.Range("J17").Value = 0.4365781710914454
$ws.Range("K17").Value = 0.135693215339233
$ws.Range("M17").Value = 0.01474926253687316
$ws.Range("N17").Value = 0.002949852507374631
$ws.Range("O17").Value = 0.05014749262536873
$ws.Range("S17").Value = 0.08849557522123894
$ws.Range("H18").Value = 0.2193877551020408
$ws.Range("I18").Value = 0.07653061224489796
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.08163265306122448
$ws.Range("M18").Value = 0.00510204081632653
$ws.Range("O18").Value = 0.04591836734693878
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.0112464854732896
$ws.Range("H19").Value = 0.1930646672914714
$ws.Range("I19").Value = 0.07029053420805999
$ws.Range("J19").Value = 0.4170571696344892
$ws.Range("K19").Value = 0.1002811621368322
$ws.Range("M19").Value = 0.01968134957825679
$ws.Range("O19").Value = 0.06279287722586692
$ws.Range("S19").Value = 0.1255857544517338